$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.828753
$ws.Range("H2").Value = 17.486259
$ws.Range("I2").Value = 0.1911291943607339
$ws.Range("J2").Value = 0.1911291943607339
$ws.Range("M2").Value = 0.01569233333333333
$ws.Range("N2").Value = 0.047077
$ws.Range("O2").Value = 0.03693539111407157
$ws.Range("P2").Value = 0.03693539111407157
$ws.Range("Q2").Value = 0.09146673499366666
$ws.Range("R2").Value = 0.823200614943
$ws.Range("S2").Value = 0.007059431547031109
$ws.Range("T2").Value = 0.007059431547031109

$ws.Range("G3").Value = 5.828753
$ws.Range("H3").Value = 17.486259
$ws.Range("I3").Value = 0.1911291943607339
$ws.Range("J3").Value = 0.1911291943607339
$ws.Range("N3").Value = 0.9690430000000001
$ws.Range("O3").Value = 0.7602859615386125
$ws.Range("P3").Value = 0.7602859615386125
$ws.Range("Q3").Value = 1.882770764459667
$ws.Range("R3").Value = 16.944936880137
$ws.Range("S3").Value = 0.1453128433126509
$ws.Range("T3").Value = 0.1453128433126509

$ws.Range("G4").Value = 5.828753
$ws.Range("H4").Value = 17.486259
$ws.Range("I4").Value = 0.1911291943607339
$ws.Range("J4").Value = 0.1911291943607339
$ws.Range("M4").Value = 0.08615233333333333
$ws.Range("N4").Value = 0.258457
$ws.Range("O4").Value = 0.202778647347316
$ws.Range("P4").Value = 0.202778647347316
$ws.Range("Q4").Value = 0.5021606713736666
$ws.Range("R4").Value = 4.519446042363
$ws.Range("S4").Value = 0.03875691950105188
$ws.Range("T4").Value = 0.03875691950105188

$ws.Range("I5").Value = 0.7732994524709527
$ws.Range("J5").Value = 0.7732994524709526
$ws.Range("M5").Value = 0.01569233333333333
$ws.Range("N5").Value = 0.047077
$ws.Range("O5").Value = 0.03693539111407157
$ws.Range("P5").Value = 0.03693539111407157
$ws.Range("Q5").Value = 0.3700699745346666
$ws.Range("R5").Value = 3.330629770812
$ws.Range("S5").Value = 0.02856211772531203
$ws.Range("T5").Value = 0.02856211772531203

$ws.Range("I6").Value = 0.7732994524709527
$ws.Range("J6").Value = 0.7732994524709526
$ws.Range("N6").Value = 0.9690430000000001
$ws.Range("O6").Value = 0.7602859615386125
$ws.Range("P6").Value = 0.7602859615386125
$ws.Range("R6").Value = 68.558392951908
$ws.Range("S6").Value = 0.5879287177791608
$ws.Range("T6").Value = 0.5879287177791608

$ws.Range("I7").Value = 0.7732994524709527
$ws.Range("J7").Value = 0.7732994524709526
$ws.Range("M7").Value = 0.08615233333333333
$ws.Range("N7").Value = 0.258457
$ws.Range("O7").Value = 0.202778647347316
$ws.Range("P7").Value = 0.202778647347316
$ws.Range("Q7").Value = 2.031717726454667
$ws.Range("R7").Value = 18.285459538092
$ws.Range("S7").Value = 0.1568086169664799
$ws.Range("T7").Value = 0.1568086169664798

$ws.Range("G8").Value = 1.084798333333333
$ws.Range("H8").Value = 3.254395
$ws.Range("I8").Value = 0.03557135316831352
$ws.Range("J8").Value = 0.03557135316831351
$ws.Range("M8").Value = 0.01569233333333333
$ws.Range("N8").Value = 0.047077
$ws.Range("O8").Value = 0.03693539111407157
$ws.Range("P8").Value = 0.03693539111407157
$ws.Range("Q8").Value = 0.01702301704611111
$ws.Range("R8").Value = 0.153207153415
$ws.Range("S8").Value = 0.001313841841728429
$ws.Range("T8").Value = 0.001313841841728428

$ws.Range("G9").Value = 1.084798333333333
$ws.Range("H9").Value = 3.254395
$ws.Range("I9").Value = 0.03557135316831352
$ws.Range("J9").Value = 0.03557135316831351
$ws.Range("N9").Value = 0.9690430000000001
$ws.Range("O9").Value = 0.7602859615386125
$ws.Range("P9").Value = 0.7602859615386125
$ws.Range("Q9").Value = 0.3504054104427778
$ws.Range("R9").Value = 3.153648693985001
$ws.Range("S9").Value = 0.02704440044680081
$ws.Range("T9").Value = 0.02704440044680081

$ws.Range("G10").Value = 1.084798333333333
$ws.Range("H10").Value = 3.254395
$ws.Range("I10").Value = 0.03557135316831352
$ws.Range("J10").Value = 0.03557135316831351
$ws.Range("M10").Value = 0.08615233333333333
$ws.Range("N10").Value = 0.258457
$ws.Range("O10").Value = 0.202778647347316
$ws.Range("P10").Value = 0.202778647347316
$ws.Range("Q10").Value = 0.09345790761277778
$ws.Range("R10").Value = 0.841121168515
$ws.Range("S10").Value = 0.007213110879784278
$ws.Range("T10").Value = 0.007213110879784277
